# Final Project Rubric - add per-student "Individual rubrik" copies.
#
# The author duplicated the "Individual rubrik" sheet three times (so each
# team member gets their own copy), filled in the newest copy with a
# specific student's name/contribution, and left the other copies (plus the
# original) untouched as blank templates.
#
# Resulting tab order (left -> right):
#   Group rubrik, Individual rubrik (4), Individual rubrik (3),
#   Individual rubrik (2), Individual rubrik

$wb = $excel.ActiveWorkbook

$original = $wb.Worksheets.Item("Individual rubrik")
$groupRubric = $wb.Worksheets.Item("Group rubrik")

# Copy 1 -> will end up as "Individual rubrik (2)" (closest to the original)
$original.Copy($null, $groupRubric)
$copy2 = $wb.Worksheets.Item(2)
$copy2.Name = "Individual rubrik (2)"

# Copy 2 -> will end up as "Individual rubrik (3)"
$original.Copy($null, $groupRubric)
$copy3 = $wb.Worksheets.Item(2)
$copy3.Name = "Individual rubrik (3)"

# Copy 3 -> will end up as "Individual rubrik (4)" (right after Group rubrik)
$original.Copy($null, $groupRubric)
$copy4 = $wb.Worksheets.Item(2)
$copy4.Name = "Individual rubrik (4)"

# Fill in the newest copy with the student's info and select the cell the
# author was last editing.
$copy4.Range("B1").Value = "Student: Kevin Zhang"
$copy4.Range("E5").Value = "Set up pygame infrastructure (game loops, game display). Coded main single player gameplay. Saved and read scores from multiple different csv files based on different settings."
$copy4.Range("E5").Select() | Out-Null
